# The workbook originally has a single sheet "Sheet1" that contains the
# sender/blacklist master data. This change turns it into a small master
# workbook: "Sheet1" is renamed to "Sender", and a second sheet "Status"
# is added holding the lookup table (value/text -> Active/Inactive) used
# to drive dropdowns elsewhere in the template.

$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet ---------------------------------------
$senderSheet = $wb.Worksheets.Item(1)
$senderSheet.Name = "Sender"

# --- Add the new "Status" sheet right after "Sender" -------------------
$statusSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $senderSheet)
$statusSheet.Name = "Status"

# --- Fill in the master lookup table -----------------------------------
$statusSheet.Range("A1").Value = "value"
$statusSheet.Range("B1").Value = "text"
$statusSheet.Range("A2").Value = 1
$statusSheet.Range("B2").Value = "Active"
$statusSheet.Range("A3").Value = 2
$statusSheet.Range("B3").Value = "Inactive"

# --- Border the whole table (thin black box around every cell) --------
$tableRange = $statusSheet.Range("A1:B3")
$tableRange.Borders.Color = 0

# --- Make the header row stand out: bold black text on a white fill ---
$headerRange = $statusSheet.Range("A1:B1")
$headerRange.Font.Bold = $true
$headerRange.Font.Color = 0
$headerRange.Interior.Color = 16777215

# --- Restore the per-sheet selections seen in the edited workbook -----
[void]$senderSheet.Range("C10").Select()
[void]$statusSheet.Range("F12").Select()

# --- "Status" is the sheet that's on top / active when reopened -------
[void]$statusSheet.Activate()
